# Fix LOM3105.xlsx "Ementa" worksheet: correct the off-by-one content
# misalignment between the label column (A) and the content columns (B/C),
# add the missing "Docentes responsaveis" (professor) rows, and fill in the
# newly-authored Portuguese objectives / summary / program / bibliography text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Portuguese text blocks (kept in variables for readability) ---
$objetivosPt = 'Apresentar ao aluno de Engenharia de Materiais ferramentas mais difundidas para análise, tratamento e exibição de dados em Engenharia. O aluno utilizará as ferramentas mais avançadas para tratamento, manipulação e apresentação de dados em planilhas eletrônicas. Ao final do curso, o estudante estará capacitado a criar correlações entre variáveis, apresentá-los em forma de gráficos e discuti-los de forma confiante. O conteúdo abordado também contempla a introdução à programação em um ambiente de fácil entendimento, de modo que o estudante compreenda e desenvolva de maneira confiante seus próprios códigos para atividades a serem desenvolvidas em outras disciplinas do curso. Os recursos aprendidos na disciplina podem ser usados em todas as demais disciplinas do curso.'
$programaResumidoPt = 'Trabalho em planilhas eletrônicas. Formatação. Importando e exportando dados. Fórmulas e Funções. Gráficos. Estatística Descritiva. Matrizes. Busca e manipulação de dados.Macros. Visual Basic for Applications. Referência a intervalos; Repetição e controle de fluxo. Programação orientada a eventos. Arrays. Classes e coleções. Introdução aos UserForms. Tabelas dinâmicas. Definindo novas funções.'
$programaPt = '- Trabalho avançado em planilhas eletrônicas:: configurando o ambiente de trabalho; Diversas variedades de planilhas; navegando por tabelas, notações de células, repetição de comandos, criação de sequências;- Formatação: verificação de condições; formatação condicional;-  Importação  e exportação de dados. Arquivos-texto, csv  (comma separated values) e outros formatos de compartilhamento eficiente de dados.- Fórmulas e Funções: Utilizando funções em planilhas eletrônicas. Usando recursos de Solver e Scenario. Sincronização de planilhas- Estatística Descritiva. Média, desvio-padrão, quartis, mediana, moda e outras características de distribuições de dados.- Matrizes: trabalhando com matrizes em planilhas eletrônicas; operações básicas: soma, multiplicação, transposição, inversão. Solução de sistemas lineares usando matrizes.- Gráficos: Gráficos de séries temporais; histogramas; gráficos de apresentação de dados (barras, setor circular, etc.) gráficos XY de correlação entre duas variáveis.- Busca e manipulação de dados: Funções de busca e identificação de dados. Operadores lógicos.- Macros: gravação, edição e utilização de sequências de comandos (macros) para automatizar tarefas. - Visual Basic for Applications: Guia desenvolvedor, VB Editor, ferramentas de depuração. Project Explorer.- Referência a intervalos: os objetos Range e Cells. Propriedades offset, resize, Columns e Rows- Repetição e controle de fluxo: laços For...Next e variações. Laços Do While/Until. Controles de fluxo If...Then...Else- Programação orientada a eventos: Níveis e parâmetros de eventos.- Arrays: declaração, arrays multidimensionais, arrays dinâmicos- Classes e coleções: Criando e usando classes, eventos de aplicação, coleções, dicionários- Introdução aos UserForms: Caixas de entrada, mensagens, botões, radio buttons;- Tabelas dinâmicas: criando e configurando uma tabela dinâmica- Definindo novas funções: estendendo as opções de funções padrão com funções definidas pelo usuário (User-Defined Functions, UDFs).'
$bibliografiaPt = '- B. JELEN, T. SYRSTAD. Excel 2016 VBA e Macros. Alta Books, 2017.- K. BLUTTMAN. Excel Fórmulas e Funções para leigos. Alta Books, 2018.- https://support.microsoft.com/pt-br/excel- https://pt-br.libreoffice.org/- http://gnumeric.org/'

# --- Content that already existed but was sitting one slot away from its label ---
$luizTadeu = '1176388 - Luiz Tadeu Fernandes Eleno'
$viktorPastoukhov = '7797767 - Viktor Pastoukhov'
$metodoContent = 'Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados. Trabalho baseado em Projeto'
$criterioContent = 'Média aritmética de trabalhos propostos ao longo do curso (60%) e do Trabalho final em grupo (40%).'
$normaContent = 'Não haverá exame de recuperação'

# 1) "Objetivos:" (row 10) was showing Prof. Eleno's name instead of the real
#    Portuguese objectives paragraph - replace it.
$ws.Range("B10:C10").Value = $objetivosPt

# 2) Insert two blank rows right after "Docentes responsaveis:" (row 12) to hold
#    the two professors, which pushes every following row down by two.
$ws.Rows.Item(13).Resize(2).Insert()

$ws.Range("B13:C13").Value = $luizTadeu
$ws.Range("B14:C14").Value = $viktorPastoukhov

# 3) "Programa resumido:" (now row 15) held a stray date - replace with the
#    real Portuguese short-syllabus paragraph.
$ws.Range("B15:C15").Value = $programaResumidoPt

# 4) "Programa:" (now row 17) held a stray professor name - replace with the
#    real, full Portuguese syllabus paragraph.
$ws.Range("B17:C17").Value = $programaPt

# 5) "Método:" (now row 20) held a stray professor name - it should hold the
#    "Aulas expositivas..." text that used to sit one row too low.
$ws.Range("B20:C20").Value = $metodoContent

# 6) "Critério:" (now row 21) should hold the grading-average text that used
#    to sit one row too low.
$ws.Range("B21:C21").Value = $criterioContent

# 7) "Norma de recuperação:" (now row 22) should hold the "não haverá exame"
#    text that used to sit one row too low.
$ws.Range("B22:C22").Value = $normaContent

# 8) "Bibliografia:" (now row 23) had no content yet - add the reading list.
$ws.Range("B23:C23").Value = $bibliografiaPt
